$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for row 2
$wsOverview.Range("G2").Value = "2016-08-23 15:22:57"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsZhCn.Range("H2").Value = "2016-08-23 15:22:52"
$wsZhCn.Range("K2").Value = "2016-08-23 15:23:25"

# de-de sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsDeDe.Range("H2").Value = "2016-08-23 15:22:57"
$wsDeDe.Range("K2").Value = "2016-08-23 15:23:33"
